# Stanford.xlsx marker-label rewrite: make "Marker:Fluorophore" labels
# self-contained/parseable instead of bare fluorophore / "Marker Fluorophore"
# text, on the "Comp controls" sheet. Also move the active-tab/selection
# focus to the "Comp controls" sheet at cell B20.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Comp controls")

$ws.Range("B4").Value  = "LIVE GREEN:FITC"
$ws.Range("B5").Value  = "CD197:PE-A"
$ws.Range("B6").Value  = "CD4:PerCP-Cy5-5-A"
$ws.Range("B7").Value  = "CD45RA:PE-Cy7"
$ws.Range("B8").Value  = "CD194:PE-Cy7"
$ws.Range("B9").Value  = "CD27:PE-Cy7"
$ws.Range("B10").Value = "CD11c:PE-Cy7"
$ws.Range("B11").Value = "CD196:PE-Cy7"
$ws.Range("B12").Value = "CD38:APC-A"
$ws.Range("B13").Value = "CD127:Alexa 647"
$ws.Range("B14").Value = "CD8:APC-H7"
$ws.Range("B15").Value = "CD45RO:APC-H7"
$ws.Range("B16").Value = "CD20:APC-Cy7-A"
$ws.Range("B17").Value = "CD3+19+20:APC-H7"
$ws.Range("B18").Value = "CD3:Pacific Blue-A"
$ws.Range("B19").Value = "HLA-DR:Am Cyan-A"

# Switch the active sheet/selection back to "Comp controls" (was "Exp
# samples"), with the cursor parked one row below the last data row.
$ws.Activate() | Out-Null
$ws.Range("B20").Select() | Out-Null
